# Sync attendance_reports: reorder "Recorded By" (column G) entries so that
# a leading "System" token is moved to the end of the comma-separated list,
# for rows whose value also references dnasr281@gmail.com or
# backup@backdoor.com (rows referencing only admin@admin.com are left as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($i = 2; $i -le $lastRow; $i++) {
    $cell = $ws.Cells.Item($i, 7)
    $val = $cell.Value2

    if ($null -eq $val) {
        continue
    }

    $parts = @($val -split ", ")

    if ($parts[0] -eq "System" -and ($val -like "*dnasr281@gmail.com*" -or $val -like "*backup@backdoor.com*")) {
        $rest = @()
        if ($parts.Length -gt 1) {
            $rest = $parts[1..($parts.Length - 1)]
        }
        $newParts = $rest + @($parts[0])
        $newVal = $newParts -join ", "
        $cell.Value = $newVal
    }
}
